# "Generate Report for Handback" - record that the zh-cn handback has come
# back in sync with en-US: populate the target/handback file columns and
# handback timestamp on the zh-cn sheet, and flip its Status.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")

# Status: "Ready for handoff" -> "Handed back: in sync with en-US"
$zh.Cells.Item(2, 3).Value = "Handed back: in sync with en-US"

# Latest Target File (I2): now points at the source markdown file, with a
# hyperlink just like the one on column A.
$zh.Cells.Item(2, 9).Value = "6d84fccc-089b-4e86-bd98-45033ca9cd74.md"
$zh.Hyperlinks.Add(
    $zh.Cells.Item(2, 9),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2969b97806e18de7e874d73a3950bef6b2e9cbea/e2e/6d84fccc-089b-4e86-bd98-45033ca9cd74.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "6d84fccc-089b-4e86-bd98-45033ca9cd74.md"
) | Out-Null
$zh.Cells.Item(2, 9).Font.Underline = $true
$zh.Cells.Item(2, 9).Font.Color = 15570276

# Latest Handback File (J2): the handed-back xliff for zh-cn.
$zh.Cells.Item(2, 10).Value = "6d84fccc-089b-4e86-bd98-45033ca9cd74.5b61347a03383bec32f73e08d8b3e3ecd619976f.zh-cn.xlf"

# Latest Handback DateTime (K2).
$zh.Cells.Item(2, 11).Value = "2016-10-10 06:17:00"

# Column autosize so the newly-populated long values are readable.
$zh.Columns.Item(3).ColumnWidth = 33
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

$ov = $wb.Worksheets.Item("Overview")
$ov.Columns.Item(5).ColumnWidth = 33
